function Set-TextValue($ws, $addr, $val) {
    # Force the cell to remain a text value: Excel would otherwise
    # auto-convert date-like / numeric-like strings (e.g. "1986-03-15",
    # "17") into real dates/numbers. Setting the format to Text first
    # keeps the literal string, then resetting the style back to Normal
    # drops the now-unneeded explicit "@" number format so the cell
    # keeps the workbook's original (unstyled) look.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
Set-TextValue $ws "A3" "Excepteur quia in ob"
Set-TextValue $ws "B3" "1986-03-15"
Set-TextValue $ws "C3" "Ahmed Mueller"
Set-TextValue $ws "D3" "Eligendi ea et tempo"
Set-TextValue $ws "E3" "2013-09-24"
Set-TextValue $ws "F3" "2003-10-14"
Set-TextValue $ws "G3" "2024-12-04"
Set-TextValue $ws "H3" "17"
Set-TextValue $ws "I3" "60"
Set-TextValue $ws "J3" "Atque numquam quia a"
Set-TextValue $ws "K3" "41"
Set-TextValue $ws "L3" "20"
Set-TextValue $ws "M3" "67"
Set-TextValue $ws "N3" "Qui esse sapiente co"
Set-TextValue $ws "O3" "Quis error exercitat"
Set-TextValue $ws "P3" "1000"
Set-TextValue $ws "Q3" "Excepturi minus non "
Set-TextValue $ws "R3" "60"
Set-TextValue $ws "S3" "99"
Set-TextValue $ws "T3" "Possimus maiores ex"
Set-TextValue $ws "U3" "55"
Set-TextValue $ws "V3" "92"

# --- Row 4 ---
Set-TextValue $ws "A4" "Excepteur quia in ob"
Set-TextValue $ws "B4" "1986-03-15"
Set-TextValue $ws "C4" "Ahmed Mueller"
Set-TextValue $ws "D4" "Eligendi ea et tempo"
Set-TextValue $ws "E4" "2013-09-24"
Set-TextValue $ws "F4" "2003-10-14"
Set-TextValue $ws "G4" "2024-12-04"
Set-TextValue $ws "H4" "17"
Set-TextValue $ws "I4" "60"
Set-TextValue $ws "J4" "Atque numquam quia a"
Set-TextValue $ws "K4" "41"
Set-TextValue $ws "L4" "20"
Set-TextValue $ws "M4" "67"
Set-TextValue $ws "N4" "Qui esse sapiente co"
Set-TextValue $ws "O4" "Quis error exercitat"
Set-TextValue $ws "P4" "1000"
Set-TextValue $ws "Q4" "Excepturi minus non "
Set-TextValue $ws "R4" "60"
Set-TextValue $ws "S4" "99"
Set-TextValue $ws "T4" "Possimus maiores ex"
Set-TextValue $ws "U4" "55"
Set-TextValue $ws "V4" "92"

# --- Row 5 ---
Set-TextValue $ws "A5" "Excepteur quia in ob"
Set-TextValue $ws "B5" "1986-03-15"
Set-TextValue $ws "C5" "Ahmed Mueller"
Set-TextValue $ws "D5" "Eligendi ea et tempo"
Set-TextValue $ws "E5" "2013-09-24"
Set-TextValue $ws "F5" "2003-10-14"
Set-TextValue $ws "G5" "2024-12-04"
Set-TextValue $ws "H5" "17"
Set-TextValue $ws "I5" "60"
Set-TextValue $ws "J5" "Atque numquam quia a"
Set-TextValue $ws "K5" "41"
Set-TextValue $ws "L5" "20"
Set-TextValue $ws "M5" "67"
Set-TextValue $ws "N5" "Qui esse sapiente co"
Set-TextValue $ws "O5" "Quis error exercitat"
Set-TextValue $ws "P5" "8"
Set-TextValue $ws "Q5" "Excepturi minus non "
Set-TextValue $ws "R5" "60"
Set-TextValue $ws "S5" "99"
Set-TextValue $ws "T5" "Possimus maiores ex"
Set-TextValue $ws "U5" "55"
Set-TextValue $ws "V5" "92"

# --- Row 6 (new row appended at the bottom of the sheet) ---
Set-TextValue $ws "A6" "Cancel"
Set-TextValue $ws "B6" "2025-01-14"
Set-TextValue $ws "C6" "2025-01-09"
Set-TextValue $ws "D6" "Rodzell Jan Gamboa Cerda"
Set-TextValue $ws "E6" "2025-01-15"
Set-TextValue $ws "F6" "2025-01-12"
Set-TextValue $ws "G6" "2025-01-15"
Set-TextValue $ws "H6" "1"
Set-TextValue $ws "I6" "2"
Set-TextValue $ws "J6" "1"
Set-TextValue $ws "K6" "2"
Set-TextValue $ws "L6" "1"
Set-TextValue $ws "M6" "1"
Set-TextValue $ws "N6" "1"
Set-TextValue $ws "O6" "2"
Set-TextValue $ws "P6" "1,10"
Set-TextValue $ws "Q6" "1"
Set-TextValue $ws "R6" "4"
Set-TextValue $ws "S6" ""
Set-TextValue $ws "T6" "3000"
Set-TextValue $ws "U6" "877"
Set-TextValue $ws "V6" "6000"

# --- Update the view's selection to match the saved state (rows 3-7 selected) ---
$ws.Rows("3:7").Select()
